# Update the date line and the division problems in the practice table.
$d = $word.ActiveDocument

# Update header date.
$d.Content.Find.Execute("2026-01-02 Friday", $true, $false, $false, $false, `
    $false, $true, 1, $false, "2026-01-03 Saturday", 2) | Out-Null

# Update the division problems. They live in rows 1, 5, 9, 13, 17 of the
# single table (5 columns each); other rows are blank spacer rows.
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Old="542÷8="; New="147÷4="},
    @{Row=1;  Col=2; Old="669÷2="; New="360÷5="},
    @{Row=1;  Col=3; Old="650÷3="; New="737÷8="},
    @{Row=1;  Col=4; Old="443÷5="; New="182÷6="},
    @{Row=1;  Col=5; Old="453÷6="; New="969÷3="},

    @{Row=5;  Col=1; Old="672÷3="; New="496÷5="},
    @{Row=5;  Col=2; Old="493÷7="; New="296÷9="},
    @{Row=5;  Col=3; Old="669÷2="; New="827÷3="},
    @{Row=5;  Col=4; Old="718÷5="; New="146÷7="},
    @{Row=5;  Col=5; Old="156÷3="; New="278÷9="},

    @{Row=9;  Col=1; Old="990÷5="; New="265÷6="},
    @{Row=9;  Col=2; Old="157÷7="; New="633÷6="},
    @{Row=9;  Col=3; Old="792÷4="; New="684÷8="},
    @{Row=9;  Col=4; Old="180÷2="; New="544÷8="},
    @{Row=9;  Col=5; Old="625÷4="; New="514÷3="},

    @{Row=13; Col=1; Old="692÷9="; New="374÷6="},
    @{Row=13; Col=2; Old="316÷3="; New="474÷2="},
    @{Row=13; Col=3; Old="453÷9="; New="641÷2="},
    @{Row=13; Col=4; Old="124÷5="; New="727÷3="},
    @{Row=13; Col=5; Old="419÷4="; New="693÷7="},

    @{Row=17; Col=1; Old="598÷5="; New="623÷6="},
    @{Row=17; Col=2; Old="517÷6="; New="678÷2="},
    @{Row=17; Col=3; Old="853÷7="; New="334÷6="},
    @{Row=17; Col=4; Old="544÷5="; New="258÷5="},
    @{Row=17; Col=5; Old="246÷2="; New="983÷2="}
)

foreach ($rep in $replacements) {
    $cellRange = $t.Cell($rep.Row, $rep.Col).Range
    $cellRange.Find.Execute($rep.Old, $true, $false, $false, $false, `
        $false, $true, 0, $false, $rep.New, 1) | Out-Null
}
